# 增加单位初始行动条位置
# Converts the ACT-bar speed / related CD fields from millisecond-scale
# integers to second-scale fractional values, and repoints the "ACT秒"
# helper column at the table's structured reference instead of the old
# literal 10000/N formula.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (field remarks/description row) -------------------------------
# ACT条速度 remark: describe the new 0-1 percent-of-bar-per-second scale.
$ws.Range("N2").Value = "施法后`n速度条`n1表示整个条`n%条/秒`n参数>0"
# 初始CD / 回转CD remarks: unit changed from milliseconds to seconds.
$ws.Range("R2").Value = "进入战斗后，首次可使用该技能的冷却CD`n(秒)"
$ws.Range("S2").Value = "首次技能使用后，每次使用技能需要间隔的最少CD(秒)"

# --- Row 6 (exported field type row) --------------------------------------
# These columns now carry fractional seconds, so the exported field type
# moves from int32 to number.
$numberCols = @("N","P","Q","R","S","U","X","Y")
foreach ($col in $numberCols) {
    $ws.Range($col + "6").Value = "number"
}

# --- Data rows 7-34 --------------------------------------------------------
# ACT条速度 (N) values were *10000 before; rebase to the 0-1 scale.
$nMap = @{5000=0.5; 10000=1; 20000=2; 2000=0.2}
# 回转CD (S) values were in milliseconds; rebase to seconds.
$sMap = @{0=0; 1000=1; 2000=2; 1500=1.5}

for ($r = 7; $r -le 34; $r++) {
    $nOld = [int]$ws.Range("N" + $r).Value()
    $ws.Range("N" + $r).Value = $nMap[$nOld]

    $sOld = [int]$ws.Range("S" + $r).Value()
    $ws.Range("S" + $r).Value = $sMap[$sOld]

    # ACT秒 = 1 / ACT条速度, expressed via the table's structured reference
    # instead of the old hard-coded 10000/N{row} formula.
    $ws.Range("O" + $r).Formula = "=1/表5[[#This Row],[ACT条速度]]"
}

# --- Sheet selection --------------------------------------------------------
$ws.Range("D6").Select()
